# Autogenerated on Fri Mar 20 2015 00:16:06 GMT+0000 (Coordinated Universal Time)
#
# In the "Belarus" MSME table, the "Enterprises density (per 1000 people)" row
# needs to be moved so that it appears immediately BEFORE the
# "Enterprises (absolute #)" row, in both the
# "Source Type: Statistical Institution" table and the
# "Source Type: SME Associations (Most Widely Used)" table. This is done by
# swapping the two whole rows (label together with its value(s)).
#
# The swap is performed via Copy/PasteSpecial (through a scratch range) so
# that the cells keep being stored as text (shared-string) values with their
# original formatting instead of being reinterpreted (e.g. as numbers) by a
# plain value re-entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")

# --- Table 1: "Source Type: Statistical Institution" ---
# Row 11 = "Enterprises (absolute #)" (A11:D11), Row 12 = "Enterprises density (per 1000 people)" (A12:D12).
# Swap them completely (label + value) so the density row comes first.
$ws.Range("A11:D11").Copy() | Out-Null
$ws.Range("Z1:AC1").PasteSpecial() | Out-Null

$ws.Range("A12:D12").Copy() | Out-Null
$ws.Range("A11:D11").PasteSpecial() | Out-Null

$ws.Range("Z1:AC1").Copy() | Out-Null
$ws.Range("A12:D12").PasteSpecial() | Out-Null

$ws.Range("Z1:AC1").Clear() | Out-Null

# --- Table 2: "Source Type: SME Associations (Most Widely Used)" ---
# Row 33 = "Enterprises (absolute #)" (A33:D33), Row 34 = "Enterprises density (per 1000 people)" (A34:D34).
# Swap them completely (label + Micro/SMEs/MSMEs values) so the density row comes first.
$ws.Range("A33:D33").Copy() | Out-Null
$ws.Range("Z3:AC3").PasteSpecial() | Out-Null

$ws.Range("A34:D34").Copy() | Out-Null
$ws.Range("A33:D33").PasteSpecial() | Out-Null

$ws.Range("Z3:AC3").Copy() | Out-Null
$ws.Range("A34:D34").PasteSpecial() | Out-Null

$ws.Range("Z3:AC3").Clear() | Out-Null

$excel.CutCopyMode = $false
